$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.955.07"
$ws.Range("E2").Value = "  +4.39%  "

# Row 3
$ws.Range("D3").Value = "2.280.49"
$ws.Range("E3").Value = "  +4.73%  "

# Row 4
$ws.Range("E4").Value = "  +0.23%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.33"
$ws.Range("E5").Value = "  +1.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.640"
$ws.Range("E6").Value = "  +3.96%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.75"
$ws.Range("E7").Value = "  +8.51%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.662"
$ws.Range("E9").Value = "  +16.15%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.30"
$ws.Range("E10").Value = "  +9.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.89"
$ws.Range("E11").Value = "  +1.76%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0967"
$ws.Range("E12").Value = "  +4.18%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.40"
$ws.Range("E13").Value = "  +7.66%  "

# Row 14
$ws.Range("E14").Value = "  +0.82%  "

# Row 15
$ws.Range("D15").Value = "2.623.95"
$ws.Range("E15").Value = "  +4.94%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.02"
$ws.Range("E16").Value = "  +4.42%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.887"
$ws.Range("E17").Value = "  +4.46%  "

# Row 18
$ws.Range("D18").Value = "2.281.02"
$ws.Range("E18").Value = "  +4.49%  "

# Row 19
$ws.Range("D19").Value = "42.926.53"
$ws.Range("E19").Value = "  +4.48%  "

# Row 20
$ws.Range("E20").Value = "  +8.10%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.33"
$ws.Range("E21").Value = "  +3.94%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.42"
$ws.Range("E22").Value = "  +2.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.73"
$ws.Range("E23").Value = "  +2.82%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  +5.32%  "

# Row 25
$ws.Range("E25").Value = "  +7.26%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("E26").Value = "  +1.70%  "

# Row 27
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.45"
$ws.Range("E28").Value = "  +1.33%  "

# Row 29
$ws.Range("E29").Value = "  -1.14%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  +9.43%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.06"
$ws.Range("E31").Value = "  +0.23%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.09"
$ws.Range("E32").Value = "  +4.42%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.52"
$ws.Range("E33").Value = "  +13.76%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.128"
$ws.Range("E34").Value = "  +5.07%  "

# Row 35
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "31.65"
$ws.Range("E35").Value = "  +29.31%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0803"
$ws.Range("E36").Value = "  +9.02%  "

# Row 37
$ws.Range("E37").Value = "  +4.34%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.51"
$ws.Range("E38").Value = "  +13.61%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.80"
$ws.Range("E39").Value = "  +6.14%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0315"
$ws.Range("E40").Value = "  +3.42%  "

# Row 41
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.45"
$ws.Range("E41").Value = "  +19.00%  "

# Row 42
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.33"
$ws.Range("E42").Value = "  +5.46%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.85"
$ws.Range("E43").Value = "  +7.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.211"
$ws.Range("E44").Value = "  +10.76%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.22"
$ws.Range("E45").Value = "  +7.97%  "

# Row 46
$ws.Range("B46").Value = "MultiversX"
$ws.Range("C46").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.37"
$ws.Range("E46").Value = "  +2.37%  "

# Row 47
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.98"
$ws.Range("E47").Value = "  -6.86%  "

# Row 48
$ws.Range("E48").Value = "  +2.71%  "

# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.19"
$ws.Range("E49").Value = "  +4.04%  "

# Row 50
$ws.Range("B50").Value = "BinanceUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.29%  "

# Row 51
$ws.Range("E51").Value = "  +4.52%  "
